$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

# The data table ("Tabela1") currently spans A1:J76 (75 data rows + header).
# Add one more data row (row 77) for 2020-05-26, growing the table to A1:J77.
$tbl = $ws.ListObjects.Item("Tabela1")
$newListRow = $tbl.ListRows.Add()
$newRowIndex = $newListRow.Range.Row

# Copy the formatting (number formats, fonts, fill, borders, alignment) of the
# previous last data row down into the freshly added row before writing the
# new values into it.
$prevRowIndex = $newRowIndex - 1
$ws.Range($ws.Cells.Item($prevRowIndex, 1), $ws.Cells.Item($prevRowIndex, 10)).Copy($ws.Range($ws.Cells.Item($newRowIndex, 1), $ws.Cells.Item($newRowIndex, 10)))

# Fill in the new row's values.
$ws.Cells.Item($newRowIndex, 1).Value = 43977
$ws.Cells.Item($newRowIndex, 2).Value = 76579
$ws.Cells.Item($newRowIndex, 3).Value = 809
$ws.Cells.Item($newRowIndex, 4).Value = 1471
$ws.Cells.Item($newRowIndex, 5).Value = 2
$ws.Cells.Item($newRowIndex, 6).Value = 8
$ws.Cells.Item($newRowIndex, 7).Value = 2
$ws.Cells.Item($newRowIndex, 8).Value = 2
$ws.Cells.Item($newRowIndex, 9).Value = 108
$ws.Cells.Item($newRowIndex, 10).Value = 0

# Match the workbook's last active selection to the newly added row, just
# like Excel leaves the selection on the row that was just entered.
$ws.Range($ws.Cells.Item($newRowIndex, 1), $ws.Cells.Item($newRowIndex, 10)).Select()
